$wb = $excel.ActiveWorkbook

# Row 17 (ALC), item id 38956
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 467.17776
$ws.Range("J17").Value = 467.17776
$ws.Range("L17").Value = 1401.53328
$ws.Range("N17").Value = -1737.53328

# Row 31 (ALC), item id 4576
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H31").Value = 1014091.8
$ws.Range("I31").Value = 1014091.8
$ws.Range("K31").Value = 3042275.4
$ws.Range("M31").Value = -3042045.4

# Row 40 (ALC), item id 5505
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 1605.8823
$ws.Range("I40").Value = 1500
$ws.Range("J40").Value = 1800
$ws.Range("K40").Value = 1500
$ws.Range("L40").Value = 1800
$ws.Range("M40").Value = -1325
$ws.Range("N40").Value = -2150

# Row 103 (ALC), item id 19909
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H103").Value = 612.8570999999999
$ws.Range("I103").Value = 470
$ws.Range("J103").Value = 670
$ws.Range("K103").Value = 1410
$ws.Range("L103").Value = 2010
$ws.Range("M103").Value = -824
$ws.Range("N103").Value = -3182

# Row 112 (ALC), item id 27960
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 2445.5
$ws.Range("J112").Value = 3191.6155
$ws.Range("L112").Value = 9574.8465
$ws.Range("N112").Value = -11790.8465

# Row 129 (ALC), item id 36115
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value = 1268.275
$ws.Range("I129").Value = 620
$ws.Range("J129").Value = 1320.8379
$ws.Range("K129").Value = 1860
$ws.Range("L129").Value = 3962.5137
$ws.Range("M129").Value = 3140
$ws.Range("N129").Value = -13962.5137

# Row 102 (ARM), item id 19945
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 4870
$ws.Range("I102").Value = 5602
$ws.Range("J102").Value = 4412.5
$ws.Range("K102").Value = 5602
$ws.Range("L102").Value = 4412.5
$ws.Range("M102").Value = -3980
$ws.Range("N102").Value = -7656.5

# Row 94 (BSM), item id 19939
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1613.25
$ws.Range("I94").Value = 900.2222
$ws.Range("J94").Value = 2530
$ws.Range("K94").Value = 900.2222
$ws.Range("L94").Value = 2530
$ws.Range("M94").Value = -449.2222
$ws.Range("N94").Value = -3432

# Row 99 (BSM), item id 19943
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 2099.12
$ws.Range("I99").Value = 1874.375
$ws.Range("J99").Value = 2498.6667
$ws.Range("K99").Value = 1874.375
$ws.Range("L99").Value = 2498.6667
$ws.Range("M99").Value = -376.375
$ws.Range("N99").Value = -5494.6667

# Row 4 (CRP), item id 3742
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 9675

# Row 16 (CRP), item id 27691
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1072.2
$ws.Range("I16").Value = 938.625
$ws.Range("J16").Value = 1606.5
$ws.Range("K16").Value = 938.625
$ws.Range("L16").Value = 1606.5
$ws.Range("M16").Value = -651.625
$ws.Range("N16").Value = -2180.5

# Row 31 (CRP), item id 44023
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2685.9546
$ws.Range("I31").Value = 1636.3684
$ws.Range("K31").Value = 1636.3684
$ws.Range("M31").Value = -1341.3684

# Row 34 (CRP), item id 44023
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 2685.9546
$ws.Range("I34").Value = 1636.3684
$ws.Range("K34").Value = 1636.3684
$ws.Range("M34").Value = -1434.3684

# Row 99 (CRP), item id 36198
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 2003.4117
$ws.Range("I99").Value = 2006.6666
$ws.Range("J99").Value = 1995.6
$ws.Range("K99").Value = 2006.6666
$ws.Range("L99").Value = 1995.6
$ws.Range("M99").Value = -508.6666
$ws.Range("N99").Value = -4991.6

# Row 105 (CRP), item id 19928
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 1859.8462
$ws.Range("I105").Value = 1982.7142
$ws.Range("J105").Value = 1716.5
$ws.Range("K105").Value = 1982.7142
$ws.Range("L105").Value = 1716.5
$ws.Range("M105").Value = -235.7141999999999
$ws.Range("N105").Value = -5210.5

# Row 107 (CRP), item id 27689
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 391.70587
$ws.Range("I107").Value = 316
$ws.Range("K107").Value = 316
$ws.Range("M107").Value = 1604

# Row 113 (CRP), item id 27691
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H113").Value = 1072.2
$ws.Range("I113").Value = 938.625
$ws.Range("J113").Value = 1606.5
$ws.Range("K113").Value = 938.625
$ws.Range("L113").Value = 1606.5
$ws.Range("M113").Value = 1231.375
$ws.Range("N113").Value = -5946.5

# Row 126 (CRP), item id 36198
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 2003.4117
$ws.Range("I126").Value = 2006.6666
$ws.Range("J126").Value = 1995.6
$ws.Range("K126").Value = 6019.9998
$ws.Range("L126").Value = 5986.799999999999
$ws.Range("M126").Value = -3549.9998
$ws.Range("N126").Value = -10926.8

# Row 4 (CUL), item id 4650
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 2966.611
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 2966.611
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = 8899.832999999999
$ws.Range("N4").Value = -9123.832999999999
$ws.Range("M4").ClearContents()

# Row 23 (CUL), item id 4858
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 306.93332
$ws.Range("I23").Value = 160.33333
$ws.Range("J23").Value = 343.58334
$ws.Range("K23").Value = 480.99999
$ws.Range("L23").Value = 1030.75002
$ws.Range("M23").Value = -245.99999
$ws.Range("N23").Value = -1500.75002

# Row 55 (CUL), item id 4733
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H55").Value = 4959
$ws.Range("I55").Value = 996
$ws.Range("J55").Value = 6657.4287
$ws.Range("K55").Value = 2988
$ws.Range("L55").Value = 19972.2861
$ws.Range("M55").Value = -2811
$ws.Range("N55").Value = -20326.2861

# Row 129 (CUL), item id 36054
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H129").Value = 3125869.5
$ws.Range("J129").Value = 3572336.5
$ws.Range("L129").Value = 10717009.5
$ws.Range("N129").Value = -10727009.5

# Row 132 (CUL), item id 43972
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 1883.3334
$ws.Range("I132").Value = 1375
$ws.Range("J132").Value = 2290
$ws.Range("K132").Value = 12375
$ws.Range("L132").Value = 20610
$ws.Range("M132").Value = -9845
$ws.Range("N132").Value = -25670

# Row 133 (CUL), item id 44073
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H133").Value = 4536.5557
$ws.Range("I133").Value = 1971.4286
$ws.Range("J133").Value = 6168.909
$ws.Range("K133").Value = 5914.2858
$ws.Range("L133").Value = 18506.727
$ws.Range("M133").Value = -854.2857999999997
$ws.Range("N133").Value = -28626.727

# Row 134 (CUL), item id 44074
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H134").Value = 3547.9644
$ws.Range("I134").Value = 2026.4706
$ws.Range("J134").Value = 5899.364
$ws.Range("K134").Value = 6079.4118
$ws.Range("L134").Value = 17698.092
$ws.Range("M134").Value = -1009.4118
$ws.Range("N134").Value = -27838.092

# Row 137 (CUL), item id 44088
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H137").Value = 2308.3333
$ws.Range("I137").Value = 1200
$ws.Range("J137").Value = 2862.5
$ws.Range("K137").Value = 3600
$ws.Range("L137").Value = 8587.5
$ws.Range("M137").Value = 1500
$ws.Range("N137").Value = -18787.5

# Row 139 (CUL), item id 44102
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H139").Value = 2208.7222
$ws.Range("I139").Value = 1576
$ws.Range("J139").Value = 2999.625
$ws.Range("K139").Value = 4728
$ws.Range("L139").Value = 8998.875
$ws.Range("M139").Value = 412
$ws.Range("N139").Value = -19278.875

# Row 93 (GSM), item id 18107
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H93").Value = 27166.666
$ws.Range("J93").Value = 27166.666
$ws.Range("L93").Value = 27166.666
$ws.Range("N93").Value = -30910.666

# Row 135 (GSM), item id 42006
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H135").Value = 51767.145
$ws.Range("J135").Value = 51767.145
$ws.Range("L135").Value = 51767.145
$ws.Range("N135").Value = -61907.145

# Row 7 (LTW), item id 36249
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1950
$ws.Range("I7").Value = 1950
$ws.Range("K7").Value = 1950
$ws.Range("M7").Value = -1838

# Row 46 (LTW), item id 5282
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 644.86206
$ws.Range("I46").Value = 750
$ws.Range("J46").Value = 628.04
$ws.Range("K46").Value = 750
$ws.Range("L46").Value = 628.04
$ws.Range("M46").Value = -562
$ws.Range("N46").Value = -1004.04

# Row 61 (LTW), item id 27740
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 15874.4375
$ws.Range("I61").Value = 23149.8
$ws.Range("J61").Value = 3748.8333
$ws.Range("K61").Value = 23149.8
$ws.Range("L61").Value = 3748.8333
$ws.Range("M61").Value = -22947.8
$ws.Range("N61").Value = -4152.8333

# Row 113 (LTW), item id 27740
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H113").Value = 15874.4375
$ws.Range("I113").Value = 23149.8
$ws.Range("J113").Value = 3748.8333
$ws.Range("K113").Value = 23149.8
$ws.Range("L113").Value = 3748.8333
$ws.Range("M113").Value = -20979.8
$ws.Range("N113").Value = -8088.8333

# Row 126 (LTW), item id 36249
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 1950
$ws.Range("I126").Value = 1950
$ws.Range("K126").Value = 5850
$ws.Range("M126").Value = -3380

# Row 2 (WVR), item id 3307
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 1999.4
$ws.Range("J2").Value = 1999.4
$ws.Range("L2").Value = 1999.4
$ws.Range("N2").Value = -2223.4

# Row 100 (WVR), item id 19981
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 9728.727999999999
$ws.Range("I100").Value = 14689.643
$ws.Range("J100").Value = 1047.125
$ws.Range("K100").Value = 29379.286
$ws.Range("L100").Value = 2094.25
$ws.Range("M100").Value = -28838.286
$ws.Range("N100").Value = -3176.25

# Row 113 (WVR), item id 27752
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 722.619
$ws.Range("I113").Value = 574.44446
$ws.Range("J113").Value = 833.75
$ws.Range("K113").Value = 1723.33338
$ws.Range("L113").Value = 2501.25
$ws.Range("M113").Value = 446.66662
$ws.Range("N113").Value = -6841.25

# Row 123 (WVR), item id 34127
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H123").Value = 33731.695
$ws.Range("J123").Value = 33731.695
$ws.Range("L123").Value = 33731.695
$ws.Range("N123").Value = -43531.695
